$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the trailing "Level 1 :" block.
#
#    The document currently ends with:
#      ... "Done" (numbered item) ...
#      <blank paragraph>
#      "Level 1 :" (bold/underlined heading)
#      <blank paragraph>
#      <numbered paragraph containing only a manual line break>
#
#    All of that, starting with the blank paragraph right after the final
#    "Done" item through to the very end of the document, is deleted.
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("Level 1 :", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$level1ParaIndex = $rng.Paragraphs.Item(1).Index

$firstParaToDelete = $d.Paragraphs.Item($level1ParaIndex - 1)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$delRange = $d.Range($firstParaToDelete.Range.Start, $lastPara.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark.
#
#    It used to sit at the very start of the "Level 0 :" paragraph; it now
#    belongs at the very end of the document, right after the text of the
#    (now last) "Done" paragraph.
# ---------------------------------------------------------------------------

$donePara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfText = $donePara.Range.End - 1

# A zero-length Range placed exactly on a paragraph-mark boundary isn't
# positioned reliably by Bookmarks.Add in this host, so a 1-character
# placeholder is inserted first, the bookmark is wrapped tightly around it,
# and the placeholder is then deleted again - leaving the (now empty)
# bookmark exactly where it belongs.
$insertPos = $d.Range($endOfText, $endOfText)
$insertPos.InsertAfter("X")

$markRange = $d.Range($endOfText, $endOfText + 1)
$d.Bookmarks.Add("_GoBack", $markRange)

$placeholder = $d.Range($endOfText, $endOfText + 1)
$placeholder.Delete()
